$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A (shifts brandName/category to B/C)
$ws.Columns.Item(1).Insert()

# New column A: Brand_id header and brand_001..brand_010 values
$ws.Range("A1").Value = "Brand_id"
for ($i = 1; $i -le 10; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = "brand_{0:D3}" -f $i
}

# Re-apply header formatting to A1 (lost on column insert) to match B1/C1
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A1").Value = "Brand_id"

# Row 11: "Lenovo" (Computers) becomes "Havells" (Mobile)
$ws.Range("B11").Value = "Havells"
$ws.Range("C11").Value = "Mobile"

# Column widths (closest achievable values to the target 24.44140625 / 26.109375 / 31.109375)
$ws.Columns.Item(1).ColumnWidth = 23.6
$ws.Columns.Item(2).ColumnWidth = 25.3
$ws.Columns.Item(3).ColumnWidth = 30.3

# View settings: zoom level and active cell selection
$ws.Select()
$excel.ActiveWindow.Zoom = 127
$ws.Range("E5").Select()
